$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column B to fit the longer project names being added
$ws.Columns.Item(2).ColumnWidth = 47.1

$ws.Range('A9').Value = 1
$ws.Range('B9').Value = 'F13a_Surgery_Assist'
$ws.Range('C9').Value = 144
$ws.Range('D9').Value = 21
$ws.Range('E9').Value = 0
$ws.Range('F9').Value = 285
$ws.Range('G9').Value = 0
$ws.Range('H9').Value = 0
$ws.Range('I9').Value = 0
$ws.Range('J9').Value = 0
$ws.Range('K9').Value = 0
$ws.Range('L9').Value = 0
$ws.Range('M9').Value = 0
$ws.Range('N9').Value = 22
$ws.Range('O9').Value = 71
$ws.Range('P9').Value = 72
$ws.Range('Q9').Value = 22
$ws.Range('R9').Value = 20
$ws.Range('S9').Value = 0.090909090909090898
$ws.Range('T9').Value = 0.05
$ws.Range('U9').Value = 2
$ws.Range('V9').Value = 1
$ws.Range('W9').Value = 2
$ws.Range('X9').Value = 2
$ws.Range('Y9').Value = 1
$ws.Range('Z9').Value = 2
$ws.Range('AA9').Value = 2
$ws.Range('AB9').Value = 6
$ws.Range('AC9').Value = 0
$ws.Range('AD9').Value = 0
$ws.Range('AE9').Value = 0
$ws.Range('AF9').Value = 72
$ws.Range('AG9').Value = 2
$ws.Range('AH9').Value = 3.2727272727272698
$ws.Range('AI9').Value = 3.22727272727272
$ws.Range('AJ9').Value = 0.27272727272727199
$ws.Range('AK9').Value = 0
$ws.Range('AL9').Value = 0
$ws.Range('AM9').Value = 0
$ws.Range('AN9').Value = 3.22727272727272
$ws.Range('AO9').Value = 0
$ws.Range('AP9').Value = 0
$ws.Range('AQ9').Value = 'NaN'
$ws.Range('AR9').Value = 1
$ws.Range('AS9').Value = 0
$ws.Range('AT9').Value = 0
$ws.Range('AU9').Value = 0
$ws.Range('AV9').Value = 0
$ws.Range('AW9').Value = 0
$ws.Range('AX9').Value = 0
$ws.Range('AY9').Value = 5
$ws.Range('AZ9').Value = 4
$ws.Range('BA9').Value = 0
$ws.Range('BB9').Value = 37
$ws.Range('BC9').Value = 4
$ws.Range('BD9').Value = 43
$ws.Range('BE9').Value = 0
$ws.Range('BF9').Value = 41
$ws.Range('BG9').Value = 101
$ws.Range('BH9').Value = 43
$ws.Range('BI9').Value = 'undefined'
$ws.Range('BJ9').Value = 'undefined'
$ws.Range('BL9').Value = 'NA'
$ws.Range('BM9').Value = 'NA'
$ws.Range('BN9').Value = 'NA'
$ws.Range('BO9').Value = 'NA'
$ws.Range('BP9').Value = 0
$ws.Range('BQ9').Value = 0
$ws.Range('BR9').Value = 21
$ws.Range('BS9').Value = 0
$ws.Range('BT9').Value = 315
$ws.Range('BU9').Value = 86
$ws.Range('BV9').Value = 86
$ws.Range('BW9').Value = 0
$ws.Range('BX9').Value = 1
$ws.Range('BY9').Value = 1
$ws.Range('BZ9').Value = 315
$ws.Range('CA9').Value = 86
$ws.Range('CB9').Value = 86
$ws.Range('CC9').Value = 1440
$ws.Range('CD9').Value = 1609
$ws.Range('CE9').Value = 1232
$ws.Range('CF9').Value = 0

$ws.Range('A10').Value = 1
$ws.Range('B10').Value = 'F13a_MedFRS_Device_Diagnostic_Software'
$ws.Range('C10').Value = 318
$ws.Range('D10').Value = 18
$ws.Range('E10').Value = 0
$ws.Range('F10').Value = 397
$ws.Range('G10').Value = 0
$ws.Range('H10').Value = 0
$ws.Range('I10').Value = 0
$ws.Range('J10').Value = 0
$ws.Range('K10').Value = 0
$ws.Range('L10').Value = 0
$ws.Range('M10').Value = 0
$ws.Range('N10').Value = 14
$ws.Range('O10').Value = 32
$ws.Range('P10').Value = 20
$ws.Range('Q10').Value = 14
$ws.Range('R10').Value = 10
$ws.Range('S10').Value = 0.42857142857142799
$ws.Range('T10').Value = 0.36363636363636298
$ws.Range('U10').Value = 4
$ws.Range('V10').Value = 4
$ws.Range('W10').Value = 4
$ws.Range('X10').Value = 4
$ws.Range('Y10').Value = 4
$ws.Range('Z10').Value = 6
$ws.Range('AA10').Value = 4
$ws.Range('AB10').Value = 3
$ws.Range('AC10').Value = 0
$ws.Range('AD10').Value = 0
$ws.Range('AE10').Value = 0
$ws.Range('AF10').Value = 20
$ws.Range('AG10').Value = 6
$ws.Range('AH10').Value = 1.4285714285714199
$ws.Range('AI10').Value = 2.2857142857142798
$ws.Range('AJ10').Value = 0.214285714285714
$ws.Range('AK10').Value = 0
$ws.Range('AL10').Value = 0
$ws.Range('AM10').Value = 0
$ws.Range('AN10').Value = 2.2857142857142798
$ws.Range('AO10').Value = 0
$ws.Range('AP10').Value = 0
$ws.Range('AQ10').Value = 'NaN'
$ws.Range('AR10').Value = 1
$ws.Range('AS10').Value = 0
$ws.Range('AT10').Value = 0
$ws.Range('AU10').Value = 0
$ws.Range('AV10').Value = 0
$ws.Range('AW10').Value = 0
$ws.Range('AX10').Value = 0
$ws.Range('AY10').Value = 29
$ws.Range('AZ10').Value = 6
$ws.Range('BA10').Value = 0
$ws.Range('BB10').Value = 27
$ws.Range('BC10').Value = 6
$ws.Range('BD10').Value = 44
$ws.Range('BE10').Value = 0
$ws.Range('BF10').Value = 33
$ws.Range('BG10').Value = 274
$ws.Range('BH10').Value = 44
$ws.Range('BI10').Value = 'undefined'
$ws.Range('BJ10').Value = 'undefined'
$ws.Range('BL10').Value = 'NA'
$ws.Range('BM10').Value = 'NA'
$ws.Range('BN10').Value = 'NA'
$ws.Range('BO10').Value = 'NA'
$ws.Range('BP10').Value = 0
$ws.Range('BQ10').Value = 0
$ws.Range('BR10').Value = 18
$ws.Range('BS10').Value = 0
$ws.Range('BT10').Value = 270
$ws.Range('BU10').Value = 88
$ws.Range('BV10').Value = 88
$ws.Range('BW10').Value = 0
$ws.Range('BX10').Value = 1
$ws.Range('BY10').Value = 1
$ws.Range('BZ10').Value = 270
$ws.Range('CA10').Value = 88
$ws.Range('CB10').Value = 88
$ws.Range('CC10').Value = 3180
$ws.Range('CD10').Value = 4229
$ws.Range('CE10').Value = 3346
$ws.Range('CF10').Value = 0

$ws.Range('A11').Value = 1
$ws.Range('B11').Value = 'S14b_E-LockBox'
$ws.Range('C11').Value = 155
$ws.Range('D11').Value = 25
$ws.Range('E11').Value = 0
$ws.Range('F11').Value = 277
$ws.Range('G11').Value = 0
$ws.Range('H11').Value = 0
$ws.Range('I11').Value = 0
$ws.Range('J11').Value = 0
$ws.Range('K11').Value = 0
$ws.Range('L11').Value = 0
$ws.Range('M11').Value = 0
$ws.Range('N11').Value = 16
$ws.Range('O11').Value = 57
$ws.Range('P11').Value = 17
$ws.Range('Q11').Value = 16
$ws.Range('R11').Value = 12
$ws.Range('S11').Value = 0.3125
$ws.Range('T11').Value = 0.0769230769230769
$ws.Range('U11').Value = 4
$ws.Range('V11').Value = 1
$ws.Range('W11').Value = 4
$ws.Range('X11').Value = 4
$ws.Range('Y11').Value = 1
$ws.Range('Z11').Value = 5
$ws.Range('AA11').Value = 4
$ws.Range('AB11').Value = 0
$ws.Range('AC11').Value = 0
$ws.Range('AD11').Value = 0
$ws.Range('AE11').Value = 0
$ws.Range('AF11').Value = 17
$ws.Range('AG11').Value = 7
$ws.Range('AH11').Value = 1.0625
$ws.Range('AI11').Value = 3.5625
$ws.Range('AJ11').Value = 0
$ws.Range('AK11').Value = 0
$ws.Range('AL11').Value = 0
$ws.Range('AM11').Value = 0
$ws.Range('AN11').Value = 3.5625
$ws.Range('AO11').Value = 0
$ws.Range('AP11').Value = 0
$ws.Range('AQ11').Value = 'NaN'
$ws.Range('AR11').Value = 1
$ws.Range('AS11').Value = 0
$ws.Range('AT11').Value = 0
$ws.Range('AU11').Value = 0
$ws.Range('AV11').Value = 0
$ws.Range('AW11').Value = 0
$ws.Range('AX11').Value = 0
$ws.Range('AY11').Value = 13
$ws.Range('AZ11').Value = 2
$ws.Range('BA11').Value = 0
$ws.Range('BB11').Value = 38
$ws.Range('BC11').Value = 2
$ws.Range('BD11').Value = 43
$ws.Range('BE11').Value = 0
$ws.Range('BF11').Value = 40
$ws.Range('BG11').Value = 112
$ws.Range('BH11').Value = 43
$ws.Range('BI11').Value = 'undefined'
$ws.Range('BJ11').Value = 'undefined'
$ws.Range('BL11').Value = 'NA'
$ws.Range('BM11').Value = 'NA'
$ws.Range('BN11').Value = 'NA'
$ws.Range('BO11').Value = 'NA'
$ws.Range('BP11').Value = 0
$ws.Range('BQ11').Value = 0
$ws.Range('BR11').Value = 25
$ws.Range('BS11').Value = 0
$ws.Range('BT11').Value = 375
$ws.Range('BU11').Value = 86
$ws.Range('BV11').Value = 86
$ws.Range('BW11').Value = 0
$ws.Range('BX11').Value = 1
$ws.Range('BY11').Value = 1
$ws.Range('BZ11').Value = 375
$ws.Range('CA11').Value = 86
$ws.Range('CB11').Value = 86
$ws.Range('CC11').Value = 1550
$ws.Range('CD11').Value = 1440
$ws.Range('CE11').Value = 1062
$ws.Range('CF11').Value = 0

$ws.Range('A12').Value = 1
$ws.Range('B12').Value = 'F13a_Yanomamo Interactive CDROM'
$ws.Range('C12').Value = 31
$ws.Range('D12').Value = 9
$ws.Range('E12').Value = 0
$ws.Range('F12').Value = 79
$ws.Range('G12').Value = 0
$ws.Range('H12').Value = 0
$ws.Range('I12').Value = 0
$ws.Range('J12').Value = 0
$ws.Range('K12').Value = 0
$ws.Range('L12').Value = 0
$ws.Range('M12').Value = 0
$ws.Range('N12').Value = 11
$ws.Range('O12').Value = 10
$ws.Range('P12').Value = 7
$ws.Range('Q12').Value = 11
$ws.Range('R12').Value = 11
$ws.Range('S12').Value = 0
$ws.Range('T12').Value = 0
$ws.Range('U12').Value = 0
$ws.Range('V12').Value = 0
$ws.Range('W12').Value = 0
$ws.Range('X12').Value = 0
$ws.Range('Y12').Value = 0
$ws.Range('Z12').Value = 0
$ws.Range('AA12').Value = 0
$ws.Range('AB12').Value = 0
$ws.Range('AC12').Value = 0
$ws.Range('AD12').Value = 0
$ws.Range('AE12').Value = 0
$ws.Range('AF12').Value = 7
$ws.Range('AG12').Value = 4
$ws.Range('AH12').Value = 0.63636363636363602
$ws.Range('AI12').Value = 0.90909090909090895
$ws.Range('AJ12').Value = 0
$ws.Range('AK12').Value = 0
$ws.Range('AL12').Value = 0
$ws.Range('AM12').Value = 0
$ws.Range('AN12').Value = 0.90909090909090895
$ws.Range('AO12').Value = 0
$ws.Range('AP12').Value = 0
$ws.Range('AQ12').Value = 'NaN'
$ws.Range('AR12').Value = 1
$ws.Range('AS12').Value = 0
$ws.Range('AT12').Value = 0
$ws.Range('AU12').Value = 0
$ws.Range('AV12').Value = 0
$ws.Range('AW12').Value = 0
$ws.Range('AX12').Value = 0
$ws.Range('AY12').Value = 2
$ws.Range('AZ12').Value = 0
$ws.Range('BA12').Value = 0
$ws.Range('BB12').Value = 14
$ws.Range('BC12').Value = 0
$ws.Range('BD12').Value = 18
$ws.Range('BE12').Value = 0
$ws.Range('BF12').Value = 14
$ws.Range('BG12').Value = 13
$ws.Range('BH12').Value = 18
$ws.Range('BI12').Value = 'undefined'
$ws.Range('BJ12').Value = 'undefined'
$ws.Range('BL12').Value = 'NA'
$ws.Range('BM12').Value = 'NA'
$ws.Range('BN12').Value = 'NA'
$ws.Range('BO12').Value = 'NA'
$ws.Range('BP12').Value = 0
$ws.Range('BQ12').Value = 0
$ws.Range('BR12').Value = 9
$ws.Range('BS12').Value = 0
$ws.Range('BT12').Value = 135
$ws.Range('BU12').Value = 36
$ws.Range('BV12').Value = 36
$ws.Range('BW12').Value = 0
$ws.Range('BX12').Value = 1
$ws.Range('BY12').Value = 1
$ws.Range('BZ12').Value = 135
$ws.Range('CA12').Value = 36
$ws.Range('CB12').Value = 36
$ws.Range('CC12').Value = 310
$ws.Range('CD12').Value = 290
$ws.Range('CE12').Value = 214
$ws.Range('CF12').Value = 0

# Select the newly added last row, matching the author's final selection state
$ws.Rows.Item(12).Select()
